# "Prototypy imprezy i skladnikow"
# Fill in the "Prototyp" (prototype) column (H) for a handful of rows on
# Arkusz1 with references to the relevant UI prototypes, and leave the
# selection on I7 (just right of the last edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("H2").Value = "p2_dodawanie_edycja_imprezy"
$ws.Range("H3").Value = "p3_wyszukanie_imprez"
$ws.Range("H4").Value = "p2_dodawanie_edycja_imprezy"
$ws.Range("H5").Value = "p4_usuniecie_imprezy"
$ws.Range("H6").Value = "p5_wyszukanie_skladnikow"
$ws.Range("H9").Value = "p6_usuniecie_skladnika"

$ws.Range("I7").Select()
